$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date value (45172 = 2023-09-03)
# for every data row (rows 2-372). Update it to the new date value
# 45175 (2023-09-06) for all those rows, preserving existing formatting.
$ws.Range("C2:C372").Value = 45175
